$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 21.13602999991758
$ws.Range("C2").Value = 381
$ws.Range("E2").Value = -0.0000007720556280982908
$ws.Range("F2").Value = 0.2377024068836954
$ws.Range("G2").Value = 3707.282554548144
$ws.Range("H2").Value = 0.5701219070552799

$ws.Range("B3").Value = 21.91985209997013
$ws.Range("C3").Value = 5
$ws.Range("F3").Value = 0.2833414185730214
$ws.Range("G3").Value = 3784.55922379831
$ws.Range("H3").Value = 0.5791916787067909

$ws.Range("B4").Value = 22.70966652997379
$ws.Range("E4").Value = 0.0000003860266280982908
$ws.Range("F4").Value = 0.302286879267628
$ws.Range("G4").Value = 3959.812384935633
$ws.Range("H4").Value = 0.5735035987151432

$ws.Range("B5").Value = 23.48373149986257
$ws.Range("F5").Value = 0.3137745999759705
$ws.Range("G5").Value = 4173.901888914051
$ws.Range("H5").Value = 0.5626325707903133

$ws.Range("B6").Value = 24.36210550014014
$ws.Range("F6").Value = 0.3200722264291359
$ws.Range("G6").Value = 4385.026548535142
$ws.Range("H6").Value = 0.5555748689430061

$ws.Range("B7").Value = 25.37542010996573
$ws.Range("D7").Value = 3
$ws.Range("E7").Value = 26.379995588
$ws.Range("F7").Value = 0.3421734567666951
$ws.Range("G7").Value = 4640.860146977853
$ws.Range("H7").Value = 0.5467826934300165

$ws.Range("B8").Value = 26.54542002996206
$ws.Range("D8").Value = 13
$ws.Range("E8").Value = 54.40294618651568
$ws.Range("F8").Value = 0.3660310163922607
$ws.Range("G8").Value = 4765.482886700025
$ws.Range("H8").Value = 0.5570352608766599

$ws.Range("B9").Value = 28.0217152499643
$ws.Range("D9").Value = 32
$ws.Range("E9").Value = 145.3991487646602
$ws.Range("F9").Value = 0.3854743807734936
$ws.Range("G9").Value = 4942.184028894884
$ws.Range("H9").Value = 0.5669905265796061

$ws.Range("B10").Value = 28.6824393299643
$ws.Range("C10").Value = 17
$ws.Range("D10").Value = 26
$ws.Range("E10").Value = 227.1082206084343
$ws.Range("F10").Value = 0.3896089912329467
$ws.Range("G10").Value = 5104.682057836884
$ws.Range("H10").Value = 0.5618849324010303

$ws.Range("B11").Value = 29.39135154996428
$ws.Range("C11").Value = 4
$ws.Range("E11").Value = 281.5750356905482
$ws.Range("F11").Value = 0.3896301042274204
$ws.Range("G11").Value = 5318.418603835355
$ws.Range("H11").Value = 0.5526332870595939

$ws.Range("B12").Value = 30.16079037998892
$ws.Range("C12").Value = 33
$ws.Range("E12").Value = 467.8542710216323
$ws.Range("F12").Value = 0.362567594658228
$ws.Range("G12").Value = 5651.598394477509
$ws.Range("H12").Value = 0.5336683230970676

$ws.Range("B13").Value = 30.95495266993477
$ws.Range("C13").Value = 11
$ws.Range("D13").Value = 61
$ws.Range("E13").Value = 508.0841732298305
$ws.Range("F13").Value = 0.326530577490232
$ws.Range("G13").Value = 5967.15084430152
$ws.Range("H13").Value = 0.518755993901109

$ws.Range("B14").Value = 31.72982493992158
$ws.Range("C14").Value = 11
$ws.Range("D14").Value = 51
$ws.Range("E14").Value = 384.2964778251018
$ws.Range("F14").Value = 0.3260579364425291
$ws.Range("G14").Value = 6228.469181448259
$ws.Range("H14").Value = 0.5094321576548875

$ws.Range("B15").Value = 31.88941800995
$ws.Range("C15").Value = 7
$ws.Range("D15").Value = 32
$ws.Range("E15").Value = 276.0190793475381
$ws.Range("F15").Value = 0.3130209330724506
$ws.Range("G15").Value = 6354.220948242343
$ws.Range("H15").Value = 0.5018619634051443

$ws.Range("B16").Value = 32.08456642994697
$ws.Range("C16").Value = 10.00000013028328
$ws.Range("D16").Value = 31
$ws.Range("E16").Value = 264.8609044457222
$ws.Range("F16").Value = 0.2874269152761762
$ws.Range("G16").Value = 6518.688981595403
$ws.Range("H16").Value = 0.4921935456735734

$ws.Range("B17").Value = 32.26831011997352
$ws.Range("C17").Value = 7
$ws.Range("D17").Value = 27
$ws.Range("E17").Value = 230.6378739437045
$ws.Range("F17").Value = 0.2768715968220207
$ws.Range("G17").Value = 6613.354510897823
$ws.Range("H17").Value = 0.4879265139468957

$ws.Range("B18").Value = 32.42996799994675
$ws.Range("E18").Value = 212.184539156656
$ws.Range("F18").Value = 0.2795966652891712
$ws.Range("G18").Value = 6785.631329923012
$ws.Range("H18").Value = 0.4779211605107154

$ws.Range("B19").Value = 32.54665038994799
$ws.Range("D19").Value = 24
$ws.Range("E19").Value = 185.1008638717538
$ws.Range("F19").Value = 0.2823518878342506
$ws.Range("G19").Value = 6873.080036417771
$ws.Range("H19").Value = 0.4735380676130058

$ws.Range("D20").Value = 17
$ws.Range("E20").Value = 132.0709356529665
$ws.Range("F20").Value = 0.2805879075827045
$ws.Range("G20").Value = 6867.761159842752
$ws.Range("H20").Value = 0.4712982121336524

$ws.Range("B21").Value = 32.17584943994828
$ws.Range("C21").Value = 4
$ws.Range("D21").Value = 15
$ws.Range("E21").Value = 108.3863609893531
$ws.Range("F21").Value = 0.2605852750073993
$ws.Range("G21").Value = 6945.25172398783
$ws.Range("H21").Value = 0.4632783766327411

$ws.Range("B22").Value = 31.94522881994698
$ws.Range("D22").Value = 8
$ws.Range("E22").Value = 37.54620682972334
$ws.Range("F22").Value = 0.2494927452418903
$ws.Range("G22").Value = 6944.706587553071
$ws.Range("H22").Value = 0.4599939308768221

$ws.Range("B23").Value = 31.6652246399495
$ws.Range("D23").Value = 2
$ws.Range("E23").Value = 0.1675222655218415
$ws.Range("F23").Value = 0.2491089332778169
$ws.Range("G23").Value = 6885.008980915103
$ws.Range("H23").Value = 0.4599155168529758

$ws.Range("B24").Value = 31.34492401994824
$ws.Range("E24").Value = 25.32281625779605
$ws.Range("F24").Value = 0.2492703839932572
$ws.Range("G24").Value = 6820.428897241818
$ws.Range("H24").Value = 0.4595740897265879

$ws.Range("B25").Value = 30.88249230004506
$ws.Range("C25").Value = 2
$ws.Range("D25").Value = 3
$ws.Range("E25").Value = 6.552910030577757
$ws.Range("F25").Value = 0.2450685097108496
$ws.Range("G25").Value = 6753.630582574115
$ws.Range("H25").Value = 0.4572724540150128

$ws.Range("B26").Value = 30.40041193995382
$ws.Range("C26").Value = 3
$ws.Range("D26").Value = 10
$ws.Range("E26").Value = 70.39519882245075
$ws.Range("F26").Value = 0.2528217540552825
$ws.Range("G26").Value = 6605.567000339182
$ws.Range("H26").Value = 0.4602241100331406

$ws.Range("B27").Value = 29.926348759952
$ws.Range("C27").Value = 11
$ws.Range("D27").Value = 21
$ws.Range("E27").Value = 136.9724303969992
$ws.Range("F27").Value = 0.2804855075434386
$ws.Range("G27").Value = 6511.630728793451
$ws.Range("H27").Value = 0.4595830139387694

$ws.Range("B28").Value = 29.46033173995674
$ws.Range("C28").Value = 13
$ws.Range("D28").Value = 35
$ws.Range("E28").Value = 186.4203940268352
$ws.Range("F28").Value = 0.2720784579469603
$ws.Range("G28").Value = 6477.261887277076
$ws.Range("H28").Value = 0.4548269354034306

$ws.Range("B29").Value = 28.99855017995705
$ws.Range("C29").Value = 28.00000071740874
$ws.Range("D29").Value = 59
$ws.Range("E29").Value = 260.0644740294942
$ws.Range("F29").Value = 0.272610607173848
$ws.Range("G29").Value = 6380.709652805436
$ws.Range("H29").Value = 0.454472178767876
